$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: append the trailing sentinel row (A522 = 0) ---
$ws1.Range("A522").Value = 0

# --- Add Sheet2 (member -> groups lookup) right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header + 107 data rows: MemberID, GroupsIDs (stringified python list).
# Many MemberID values are long all-digit strings; writing them straight to
# .Value would get auto-coerced to numbers by Excel. Instead we stash every
# value behind a literal-text formula (="...") and then Copy/PasteSpecial
# (values-only) over the same range, which flips the cells to plain shared-
# string literals without touching any cell formatting/styles.
$data = @(
    ,@('MemberID', 'GroupsIDs')
    ,@('667187513392561', '[''324309874271040'', ''255488514638473'', ''334271300068285'']')
    ,@('774295729293068', '[''324309874271040'', ''255488514638473'']')
    ,@('621182804659975', '[''206494919465453'', ''255488514638473'']')
    ,@('407161542767385', '[''324309874271040'', ''255488514638473'']')
    ,@('691736924273562', '[''255488514638473'', ''334271300068285'']')
    ,@('633931720052579', '[''324309874271040'', ''255488514638473'']')
    ,@('1616975745197451', '[''324309874271040'', ''255488514638473'']')
    ,@('482088201930268', '[''324309874271040'', ''255488514638473'']')
    ,@('336335959888767', '[''255488514638473'', ''334271300068285'']')
    ,@('1554985268075193', '[''324309874271040'', ''255488514638473'']')
    ,@('814428541931662', '[''324309874271040'', ''255488514638473'']')
    ,@('501399673335368', '[''324309874271040'', ''255488514638473'']')
    ,@('1585793238309987', '[''324309874271040'', ''255488514638473'']')
    ,@('549638481834556', '[''324309874271040'', ''255488514638473'']')
    ,@('876041585781332', '[''324309874271040'', ''255488514638473'', ''334271300068285'']')
    ,@('771972242840206', '[''324309874271040'', ''255488514638473'']')
    ,@('640052912772312', '[''324309874271040'', ''255488514638473'']')
    ,@('400671996754563', '[''324309874271040'', ''255488514638473'']')
    ,@('793081574090815', '[''324309874271040'', ''255488514638473'']')
    ,@('1000941419935908', '[''324309874271040'', ''255488514638473'']')
    ,@('1542466826000442', '[''324309874271040'', ''255488514638473'']')
    ,@('808510189214327', '[''324309874271040'', ''255488514638473'']')
    ,@('680412775390168', '[''324309874271040'', ''255488514638473'']')
    ,@('569530223183696', '[''324309874271040'', ''255488514638473'']')
    ,@('876704639029711', '[''324309874271040'', ''334271300068285'']')
    ,@('755493964529658', '[''324309874271040'', ''255488514638473'', ''334271300068285'']')
    ,@('484785774992731', '[''324309874271040'', ''255488514638473'']')
    ,@('599462043488316', '[''324309874271040'', ''255488514638473'']')
    ,@('728792203837815', '[''324309874271040'', ''255488514638473'']')
    ,@('884743478233301', '[''324309874271040'', ''255488514638473'']')
    ,@('10204419721387878', '[''324309874271040'', ''255488514638473'']')
    ,@('725408254223667', '[''324309874271040'', ''255488514638473'']')
    ,@('787357154676765', '[''324309874271040'', ''255488514638473'']')
    ,@('857964870890641', '[''324309874271040'', ''255488514638473'']')
    ,@('963700236991334', '[''324309874271040'', ''255488514638473'']')
    ,@('820661414659110', '[''324309874271040'', ''255488514638473'']')
    ,@('689786061134490', '[''324309874271040'', ''206494919465453'']')
    ,@('858101107537227', '[''324309874271040'', ''255488514638473'']')
    ,@('1578393622381080', '[''324309874271040'', ''255488514638473'']')
    ,@('814692308592705', '[''255488514638473'', ''334271300068285'']')
    ,@('640790446046729', '[''324309874271040'', ''255488514638473'']')
    ,@('806424042751527', '[''324309874271040'', ''255488514638473'']')
    ,@('1567026680207269', '[''324309874271040'', ''255488514638473'', ''334271300068285'']')
    ,@('315557471977295', '[''324309874271040'', ''255488514638473'']')
    ,@('1527519694201535', '[''324309874271040'', ''255488514638473'']')
    ,@('548810601889073', '[''324309874271040'', ''255488514638473'']')
    ,@('1530768973863166', '[''255488514638473'', ''334271300068285'']')
    ,@('1547593428818706', '[''324309874271040'', ''255488514638473'']')
    ,@('616304628475710', '[''324309874271040'', ''255488514638473'']')
    ,@('313504398857318', '[''324309874271040'', ''255488514638473'']')
    ,@('678794352241097', '[''324309874271040'', ''255488514638473'']')
    ,@('851976011533451', '[''324309874271040'', ''255488514638473'', ''334271300068285'']')
    ,@('566915410111458', '[''324309874271040'', ''255488514638473'']')
    ,@('645377155584155', '[''324309874271040'', ''255488514638473'']')
    ,@('929705337042553', '[''324309874271040'', ''255488514638473'']')
    ,@('778111352243342', '[''324309874271040'', ''255488514638473'']')
    ,@('564297323701188', '[''324309874271040'', ''255488514638473'']')
    ,@('641351659309361', '[''324309874271040'', ''255488514638473'']')
    ,@('742870925794424', '[''324309874271040'', ''255488514638473'']')
    ,@('825063557536678', '[''324309874271040'', ''255488514638473'', ''334271300068285'']')
    ,@('1555238544693890', '[''324309874271040'', ''255488514638473'']')
    ,@('1016370118380238', '[''324309874271040'', ''255488514638473'']')
    ,@('839048342805816', '[''324309874271040'', ''255488514638473'', ''334271300068285'']')
    ,@('717272155053183', '[''324309874271040'', ''255488514638473'']')
    ,@('551193358328556', '[''324309874271040'', ''255488514638473'', ''334271300068285'']')
    ,@('1537405156507717', '[''324309874271040'', ''206494919465453'']')
    ,@('798334086880359', '[''324309874271040'', ''255488514638473'']')
    ,@('591137314354302', '[''324309874271040'', ''255488514638473'']')
    ,@('877708612262360', '[''324309874271040'', ''255488514638473'']')
    ,@('10203220024557016', '[''324309874271040'', ''255488514638473'']')
    ,@('710461015679951', '[''324309874271040'', ''255488514638473'']')
    ,@('397932313693174', '[''324309874271040'', ''255488514638473'']')
    ,@('865923580113918', '[''324309874271040'', ''255488514638473'']')
    ,@('1546713195576368', '[''206494919465453'', ''255488514638473'']')
    ,@('748789238538973', '[''324309874271040'', ''255488514638473'']')
    ,@('514255645384204', '[''324309874271040'', ''255488514638473'']')
    ,@('1679238515636079', '[''324309874271040'', ''255488514638473'']')
    ,@('753306174747351', '[''324309874271040'', ''255488514638473'']')
    ,@('538192309656914', '[''324309874271040'', ''255488514638473'', ''334271300068285'']')
    ,@('895705570470293', '[''324309874271040'', ''255488514638473'']')
    ,@('579563348812799', '[''324309874271040'', ''255488514638473'']')
    ,@('625591577570435', '[''324309874271040'', ''255488514638473'']')
    ,@('767213676691928', '[''255488514638473'', ''334271300068285'']')
    ,@('1516260108639780', '[''324309874271040'', ''255488514638473'', ''334271300068285'']')
    ,@('1501169373481444', '[''324309874271040'', ''255488514638473'', ''334271300068285'']')
    ,@('849156381809820', '[''324309874271040'', ''255488514638473'']')
    ,@('825300097508778', '[''324309874271040'', ''255488514638473'']')
    ,@('10203649275811944', '[''324309874271040'', ''255488514638473'']')
    ,@('10203566138286951', '[''324309874271040'', ''255488514638473'']')
    ,@('541962275906605', '[''324309874271040'', ''255488514638473'']')
    ,@('309341389262137', '[''324309874271040'', ''255488514638473'']')
    ,@('832692130124011', '[''324309874271040'', ''255488514638473'']')
    ,@('900715926619149', '[''324309874271040'', ''255488514638473'']')
    ,@('1486838111528993', '[''324309874271040'', ''255488514638473'']')
    ,@('782035871871854', '[''324309874271040'', ''255488514638473'']')
    ,@('823200181077355', '[''324309874271040'', ''255488514638473'']')
    ,@('777204325679347', '[''324309874271040'', ''255488514638473'']')
    ,@('10152692769938915', '[''324309874271040'', ''255488514638473'']')
    ,@('908926152475002', '[''206494919465453'', ''255488514638473'']')
    ,@('844212072308634', '[''324309874271040'', ''206494919465453'', ''255488514638473'', ''334271300068285'']')
    ,@('1601618586738184', '[''324309874271040'', ''206494919465453'', ''255488514638473'']')
    ,@('1572362456311048', '[''324309874271040'', ''255488514638473'', ''334271300068285'']')
    ,@('10205541649171250', '[''324309874271040'', ''255488514638473'']')
    ,@('1550815785165530', '[''324309874271040'', ''255488514638473'']')
    ,@('534532093316239', '[''324309874271040'', ''255488514638473'', ''334271300068285'']')
    ,@('1039052362786925', '[''324309874271040'', ''255488514638473'']')
    ,@('384146325094022', '[''206494919465453'', ''255488514638473'']')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $ws2.Cells.Item($r, 1).Formula = '="' + $data[$i][0] + '"'
    $ws2.Cells.Item($r, 2).Formula = '="' + $data[$i][1] + '"'
}

# Convert the literal-text formulas to static values in one shot (keeps
# xl/styles.xml untouched -- no new number formats / cell styles created).
$used = $ws2.Range("A1:B" + $data.Count)
$used.Copy()
$used.PasteSpecial(-4163) | Out-Null

# Restore Sheet1 as the active/selected sheet (matches the original state).
$ws1.Activate()
